$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = 'ام العبد الترتوري 2'
$ws.Range("A30").Value = 'ام سائد شركة فيجين'
$ws.Range("A38").Value = 'عبد اللطيف عبيد'
$ws.Range("A52").Value = 'ام محمود مهنا'
$ws.Range("A56").Value = 'محمد فوره'
$ws.Range("A58").Value = 'سعدو خلف'
$ws.Range("A76").Value = 'وفيق2'
$ws.Range("A85").Value = 'مصطفى بشير'
$ws.Range("A89").Value = 'ابويونس شملخ'
$ws.Range("A100").Value = 'محمد شعت'
$ws.Range("A102").ClearContents()
$ws.Range("A103").Value = 'ابوفياض'
$ws.Range("A109").Value = 'ابومدلله سمارت فون'
$ws.Range("A111").Value = 'ابوضياء الخالدي'
$ws.Range("A113").Value = 'عمار ابوضاهر'
$ws.Range("A115").Value = 'الدانا'
$ws.Range("A119").Value = 'فراس دير البلح'
$ws.Range("A124").Value = 'اميره الواديه'
$ws.Range("A125").Value = 'شادي ابوحصيره'
$ws.Range("A126").Value = 'ابوحمده'
$ws.Range("A129").Value = 'حسن العبادله'
$ws.Range("A130").Value = 'خالد محمود'
$ws.Range("A131").Value = 'الدكتور القيشاوي'
$ws.Range("A132").Value = 'محمد الخزندار'
$ws.Range("A133").Value = 'حميد'

$ws.Range("A134:A138").EntireRow.Delete()
